$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

# Pre-format the Price/Volume data range as Text so numeric-looking strings
# (e.g. "1.002", "0.3660") are stored verbatim instead of being re-parsed as numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '27.505.62'
$ws.Range('E2').Value = '  -0.80%  '
$ws.Range('D3').Value = '1.831.07'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '312.67'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Value = '0.4287'
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('D8').Value = '0.3660'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.07269'
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('D10').Value = '0.8658'
$ws.Range('E10').Value = '  -1.49%  '
$ws.Range('D11').Value = '20.63'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').Value = '1.865.79'
$ws.Range('E12').Value = '  +3.22%  '
$ws.Range('D13').Value = '5.399'
$ws.Range('E13').Value = '  +0.76%  '
$ws.Range('D14').Value = '6.543'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = '80.60'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').Value = '0.000008907'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = '15.43'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('D21').Value = '27.674.38'
$ws.Range('E21').Value = '  +0.30%  '
$ws.Range('D22').Value = '5.150'
$ws.Range('E22').Value = '  +3.49%  '
$ws.Range('D23').Value = '10.87'
$ws.Range('E23').Value = '  +5.20%  '
$ws.Range('D24').Value = '2.109.28'
$ws.Range('E24').Value = '  +2.23%  '
$ws.Range('E25').Value = '  -0.31%  '
$ws.Range('D26').Value = '154.58'
$ws.Range('E26').Value = '  -1.05%  '
$ws.Range('D27').Value = '18.87'
$ws.Range('E27').Value = '  +1.29%  '
$ws.Range('D28').Value = '5.123'
$ws.Range('E28').Value = '  -2.41%  '
$ws.Range('D29').Value = '114.11'
$ws.Range('E29').Value = '  -5.03%  '
$ws.Range('D30').Value = '1.836'
$ws.Range('E30').Value = '  -2.62%  '
$ws.Range('D31').Value = '0.08850'
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').Value = '3.010'
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '0.7521'
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('D34').Value = '4.543'
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('D35').Value = '1.133'
$ws.Range('E35').Value = '  +0.45%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('D37').Value = '1.088'
$ws.Range('E37').Value = '  -1.96%  '
$ws.Range('D38').Value = '0.05324'
$ws.Range('E38').Value = '  -2.20%  '
$ws.Range('D39').Value = '0.01938'
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').Value = '2.793'
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('D41').Value = '0.5085'
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('D42').Value = '0.1662'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').Value = '6.562'
$ws.Range('E43').Value = '  -0.76%  '
$ws.Range('D44').Value = '8.327'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('D45').Value = '10.39'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('D46').Value = '105.96'
$ws.Range('E46').Value = '  +0.38%  '
$ws.Range('D47').Value = '0.06499'
$ws.Range('E47').Value = '  -0.87%  '
$ws.Range('D48').Value = '0.4684'
$ws.Range('E48').Value = '  +0.37%  '
$ws.Range('D50').Value = '1.616'
$ws.Range('E50').Value = '  -1.45%  '
$ws.Range('D51').Value = '63.95'
$ws.Range('E51').Value = '  -0.81%  '

# Restore the default (Normal) style so no stray NumberFormat is left on the cells,
# matching the original workbook formatting.
$dataRange.Style = "Normal"

